# Update Name of Algo
# Apply updated RandomForest imputation results to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.687599999999994
$ws.Range("C3").Value = -11.32469999999999
$ws.Range("A4").Value = -21.64730000000002
$ws.Range("B4").Value = 4.830800000000002
$ws.Range("C4").Value = -11.40029999999999
$ws.Range("E4").Value = 13.60350000000002
$ws.Range("B5").Value = 5.4781
$ws.Range("D5").Value = -9.205699999999988
$ws.Range("A6").Value = -21.56940000000001
$ws.Range("B6").Value = 5.783599999999998
$ws.Range("A7").Value = -21.45250000000001
$ws.Range("A8").Value = -21.42310000000001
$ws.Range("B8").Value = 4.798999999999996
$ws.Range("C9").Value = -11.9018
$ws.Range("C11").Value = -14.22100000000001
$ws.Range("E12").Value = 11.02890000000001
$ws.Range("C14").Value = -11.72519999999999
$ws.Range("A16").Value = -21.23570000000002
$ws.Range("B16").Value = 4.959499999999997
$ws.Range("E16").Value = 13.08820000000001
$ws.Range("E17").Value = 13.04060000000001
$ws.Range("C18").Value = -14.45300000000002
$ws.Range("A20").Value = -22.59730000000001
$ws.Range("D20").Value = -8.214200000000003
$ws.Range("E20").Value = 13.26239999999999
$ws.Range("A21").Value = -20.60049999999999
$ws.Range("B22").Value = 5.3669
$ws.Range("C25").Value = -11.55539999999999
$ws.Range("E25").Value = 13.61570000000001
